$d = $word.ActiveDocument

# --- 1. Add a comment on the word "Year" in the table header row ---
$word.UserName = "Gabriella Meltzer"
$word.UserInitials = "GM"

$rng = $d.Content
$found = $rng.Find.Execute("Year", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $d.Comments.Add($rng, "Insert 2008") | Out-Null
}

# --- 2. Delete the last table row (2018 / Alberto.../354) ---
$tbl = $d.Tables.Item(1)
$tbl.Rows.Item($tbl.Rows.Count).Delete()
